# Update Number_of_Inclusions (col B) and recompute
# Number_of_Inclusions_per_Nucleus (col D = B / C) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new value for column B (Number_of_Inclusions)
$updates = @{
    11 = 12
    12 = 2
    14 = 8
    22 = 4
    23 = 25
    29 = 1
    31 = 2
    32 = 5
    53 = 2
}

foreach ($row in $updates.Keys) {
    $newB = $updates[$row]

    $ws.Cells.Item($row, 2).Value2 = $newB

    $c = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value2 = $newB / $c
}
